# Updated user stories to be more precise.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 3 (story #1, Film Producer) ---
$ws.Range("D3").Value = "Filter the project proposals by different categories (expected budget, length, genres, etc.)"
$ws.Range("E3").Value = "Easily find the types of project my company is interested in."

# --- Row 4 (story #2, Film Producer) ---
$ws.Range("B4").Formula = "=B3+1"
$ws.Range("D4").Value = "Obtain the detailed information about the producer (years of experience, production company, studies, etc.)"
$ws.Range("E4").Value = "Reach out the them if I am interested in a potential collaboration."

# --- Row 5 (story #3, Film Producer) ---
$ws.Range("B5").Formula = "=B4+1"
$ws.Range("D5").Value = "Share my proposal while maintaining my rights to it"
$ws.Range("E5").Value = "Be contacted by other producers interested in it."

# --- Row 6 (story #4, Film Producer) ---
$ws.Rows.Item(6).RowHeight = 45.5
$ws.Range("B6").Formula = "=B5+1"
$ws.Range("D6").Value = "Remove any of my proposals from the forum"
$ws.Range("E6").Value = "Stop being contacted about it when my company has already decided which other production will be collaborating with us."

# --- Row 7 (story #5, Film Producer) ---
$ws.Rows.Item(7).RowHeight = 55.5
$ws.Range("B7").Formula = "=B6+1"
$ws.Range("C7").Value = "Film Producer"
$ws.Range("D7").Value = "Customize the graphs (e.g. focus on a specific period, remove certain genres or movies from distributors that I am not interested in, have the revenue in a logarithmic scale, increase the font size of the labels, etc.)"
$ws.Range("E7").Value = "Have my graphs laid out in a way that facilitates the research that I want to conduct."

# --- Row 8 (story #6, Film Producer) ---
$ws.Rows.Item(8).RowHeight = 56.5
$ws.Range("B8").Formula = "=B7+1"
$ws.Range("C8").Value = "Film Producer"
$ws.Range("D8").Value = "View different types of graphs (pie charts, bar and linear plots) that display the correlation between different variables (i.e. revenue, genre, ratings, distribution companies, runtime)"
$ws.Range("E8").Value = "Have enough information to carry out my desired research."

# --- Row 9 (story #7, Film Producer) - previously blank ---
$ws.Range("B9").Formula = "=B8+1"
$ws.Range("C9").Value = "Film Producer"
$ws.Range("D9").Value = "Save/Download graphs in various formats (png, jpg, jpeg, tif, etc.)"
$ws.Range("E9").Value = "Include them in the research paper that I will be producing, or show it to my boss to enage in discussion."

# --- Row 10 (story #8, Film Producer) - previously blank ---
$ws.Range("B10").Formula = "=B9+1"
$ws.Range("C10").Value = "Film Producer"
$ws.Range("D10").Value = "Be able to contact the customer support team (either via email or a phone call)"
$ws.Range("E10").Value = "Ask any questions I may have regarding the use, performance or data protection of the App."

# --- Row 11 (story #9, Film Producer) - previously blank ---
$ws.Range("B11").Formula = "=B10+1"
$ws.Range("C11").Value = "Film Producer"
$ws.Range("D11").Value = "Have the option to either freeze or delete my account"
$ws.Range("E11").Value = "Stop using the app for a period of time or indefinitely. "

# --- Row 12 (story #10, Advertiser) - previously blank ---
$ws.Range("B12").Formula = "=B11+1"
$ws.Range("C12").Value = "Advertiser"
$ws.Range("D12").Value = "Ensure that my add slot is visible and displayed regularly (in case there are multiple ads being shown by the app)"
$ws.Range("E12").Value = "Maximize the number of users that see our promotion."

# --- Row 13 (story #11, Passive Investor) - previously blank ---
$ws.Range("B13").Formula = "=B12+1"
$ws.Range("C13").Value = "Passive Investor"
$ws.Range("D13").Value = "Buy a small amount of company stock"
$ws.Range("E13").Value = "Gain money if the company experiences financial growth."

# --- Row 14 (story #12, Active Investor) - previously blank ---
$ws.Rows.Item(14).RowHeight = 41.5
$ws.Range("B14").Formula = "=B13+1"
$ws.Range("C14").Value = "Active Investor"
$ws.Range("D14").Value = "But a significant amount of company stock"
$ws.Range("E14").Value = "Negotiate a position in the director's board, have the right to vote on important business decisions and be informed about them, etc."

# --- View state: scroll the frozen pane down and select D12 ---
$ws.Range("A8").Select()
$ws.Range("D12").Select()
